$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A33").Value = "19 marras"
$ws.Range("B33").Value = "16.30-17.30"
$ws.Range("C33").Value = "Rajaavasta rakenteesta poistaminen. Nopea johdanto partiotioiviin rakenteisiin."

$ws.Range("B33").NumberFormat = "h:mm"
$ws.Range("C33").WrapText = $true

$ws.Rows.Item(33).RowHeight = 43.5

$ws.Range("C33").Select() | Out-Null
